# certidao_obito.docx template update
#  - merge the "nota_/bg" placeholder runs into a single "nota_bg" run
#  - change "DRH-1/202" to "DP-1/202" and move the hidden _GoBack bookmark
#    so it now sits right after "DP" (instead of at the very end of the doc)
#  - reword "Diretoria de Recursos Humanos" -> "Diretoria de Pessoal" in the
#    "Foi encaminhada..." paragraph, splitting it into three runs
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "nota_" + "bg" -> single run "nota_bg"
#    (also drops the now-stale w:proofErr "gramStart" marker between them)
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("nota_bg", $false, $false, $false, $false, $false, $true, 1, $false, "nota_bg", 2)

# ---------------------------------------------------------------------
# 2) "DRH-1/202" -> "DP-1/202", with the _GoBack bookmark ending up right
#    after "DP". Re-adding a bookmark with an existing name relocates it,
#    so the stale bookmark at the end of the document disappears on its
#    own once this one is placed.
# ---------------------------------------------------------------------
$text = $d.Content.Text
$idx = $text.IndexOf("DRH-1/202")
$bmPos = $idx + 3
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Replace "DRH" with "DP", spanning the preceding "/" so the stale
# w:proofErr "gramEnd" marker (sitting right before "DRH") is cleaned up
# too. This merges "/" into the same run as "DP-1/2025." - immediately
# re-split "/" back off with a scratch bookmark (added then deleted)
# so it stays its own run, matching the original structure.
$text = $d.Content.Text
$idx = $text.IndexOf("DRH")
$r = $d.Range($idx - 1, $idx + 3)
$r.Text = "/DP"

$text = $d.Content.Text
$idx = $text.IndexOf("/DP")
$splitPos = $idx + 1
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_scratch_split", $splitRange)
$d.Bookmarks.Item("_scratch_split").Delete()

# ---------------------------------------------------------------------
# 3) "Foi encaminhada à Diretoria de Recursos Humanos do CBMAM, cópia da"
#    -> "Foi encaminhada" / " à Diretoria de Pessoal" / " do CBMAM, cópia da"
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("Recursos Humanos", $false, $false, $false, $false, $false, $true, 1, $false, "Pessoal", 2)

$text = $d.Content.Text
$idx = $text.IndexOf("Foi encaminhada")
$splitPos1 = $idx + "Foi encaminhada".Length
$splitPos2 = $idx + "Foi encaminhada à Diretoria de Pessoal".Length

$r1 = $d.Range($splitPos1, $splitPos1)
$d.Bookmarks.Add("_scratch_split", $r1)
$d.Bookmarks.Item("_scratch_split").Delete()

$r2 = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("_scratch_split", $r2)
$d.Bookmarks.Item("_scratch_split").Delete()

# The text edit above re-normalises every run in this paragraph, which
# also silently glues the unrelated "CERTIDÃO" / " DE ÓBITO" runs (further
# along in the same paragraph) back together. Re-split them so that part
# of the paragraph stays exactly as it was before.
$text = $d.Content.Text
$idx = $text.IndexOf("CERTIDÃO DE ÓBITO")
$splitPos3 = $idx + "CERTIDÃO".Length
$r3 = $d.Range($splitPos3, $splitPos3)
$d.Bookmarks.Add("_scratch_split", $r3)
$d.Bookmarks.Item("_scratch_split").Delete()

Write-Output "done"
